$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170101881027222
$ws.Range("B1").Value = 2.439114332199097
$ws.Range("D1").Value = 2.367249488830566
$ws.Range("E1").Value = 1.234473347663879
